$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain as
# literal text (matching the original "General"-formatted text cells,
# e.g. keeping trailing zeros / exact decimal text). Force text format,
# assign the value, then reset the style so no stray number format is
# left attached to the cell.
$forceTextAddrs = @(
    "D5",
    "D6",
    "D8",
    "D10",
    "D13",
    "D15",
    "D18",
    "D19",
    "D21",
    "D24",
    "D25",
    "D26",
    "D27",
    "D31",
    "D36",
    "D38",
    "D39",
    "D40",
    "D42",
    "D45",
    "D50"
)
foreach ($addr in $forceTextAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the forced-text numeric-looking values
$ws.Range("D5").Value = "227.90"
$ws.Range("D6").Value = "0.603"
$ws.Range("D8").Value = "36.19"
$ws.Range("D10").Value = "0.0695"
$ws.Range("D13").Value = "11.29"
$ws.Range("D15").Value = "0.646"
$ws.Range("D18").Value = "70.04"
$ws.Range("D19").Value = "245.39"
$ws.Range("D21").Value = "11.48"
$ws.Range("D24").Value = "2.27"
$ws.Range("D25").Value = "171.20"
$ws.Range("D26").Value = "8.11"
$ws.Range("D27").Value = "17.52"
$ws.Range("D31").Value = "1.25"
$ws.Range("D36").Value = "0.663"
$ws.Range("D38").Value = "0.0188"
$ws.Range("D39").Value = "2.34"
$ws.Range("D40").Value = "82.38"
$ws.Range("D42").Value = "0.950"
$ws.Range("D45").Value = "13.48"
$ws.Range("D50").Value = "103.37"

foreach ($addr in $forceTextAddrs) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining text updates (names, links, price strings that are not valid
# numbers, and percentage-change strings) can be assigned directly since
# Excel will keep them as text.
$ws.Range("D2").Value = "34.422.13"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.806.67"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  +6.24%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +3.34%  "
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "2.066.43"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").Value = "1.822.38"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "34.413.52"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "0.0₃0787"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +8.61%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  +7.79%  "
$ws.Range("E27").Value = "  +4.45%  "
$ws.Range("E28").Value = "  +3.96%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("E34").Value = "  -3.04%  "
$ws.Range("D35").Value = "1.382.14"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E39").Value = "  -11.02%  "
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  +6.81%  "
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("D48").Value = "1.967.73"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").Value = "0.0₆0126"
$ws.Range("E51").Value = "  -4.10%  "
